# Apply updated cryptos list values (price + 1h volume change)
# GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.557.81'
$ws.Range("E2").Value = '  -1.45%  '
$ws.Range("D3").Value = '3.033.35'
$ws.Range("E3").Value = '  -4.43%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''575.42'
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("D6").Value = '''129.44'
$ws.Range("E6").Value = '  -4.14%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.031.92'
$ws.Range("E8").Value = '  -4.40%  '
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("E10").Value = '  -3.21%  '
$ws.Range("D11").Value = '''5.22'
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  -4.33%  '
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("D14").Value = '''33.41'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").Value = '3.536.24'
$ws.Range("D17").Value = '61.691.82'
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("D18").Value = '3.035.77'
$ws.Range("E18").Value = '  -4.33%  '
$ws.Range("E19").Value = '  -3.05%  '
$ws.Range("D20").Value = '''442.57'
$ws.Range("E20").Value = '  -2.68%  '
$ws.Range("D21").Value = '''13.34'
$ws.Range("E21").Value = '  -4.33%  '
$ws.Range("D22").Value = '''0.668'
$ws.Range("E22").Value = '  -4.75%  '
$ws.Range("E23").Value = '  -4.64%  '
$ws.Range("D24").Value = '''80.34'
$ws.Range("E24").Value = '  -3.81%  '
$ws.Range("D25").Value = '''12.70'
$ws.Range("E25").Value = '  -4.21%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '''0.999'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -5.41%  '
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = '''7.41'
$ws.Range("E30").Value = '  -4.08%  '
$ws.Range("E31").Value = '  -6.98%  '
$ws.Range("D32").Value = '''25.70'
$ws.Range("E32").Value = '  -5.66%  '
$ws.Range("D33").Value = '''0.0959'
$ws.Range("E33").Value = '  -6.71%  '
$ws.Range("D34").Value = '''2.31'
$ws.Range("E34").Value = '  -2.83%  '
$ws.Range("D35").Value = '''0.968'
$ws.Range("E35").Value = '  -6.15%  '
$ws.Range("D36").Value = '''5.64'
$ws.Range("E36").Value = '  -4.58%  '
$ws.Range("D37").Value = '''50.30'
$ws.Range("E37").Value = '  -1.57%  '
$ws.Range("D38").Value = '0.0₃0696'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("E39").Value = '  -3.55%  '
$ws.Range("D40").Value = '''7.86'
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("E42").Value = '  -7.52%  '
$ws.Range("D43").Value = '''375.63'
$ws.Range("E43").Value = '  -5.06%  '
$ws.Range("D44").Value = '2.672.08'
$ws.Range("D46").Value = '''123.55'
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("E47").Value = '  -4.89%  '
$ws.Range("D48").Value = '''33.90'
$ws.Range("E48").Value = '  -4.66%  '
$ws.Range("E49").Value = '  -6.56%  '
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("D51").Value = '''23.59'
$ws.Range("E51").Value = '  -6.56%  '
